$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is plain-text (coin names, URLs, and percent-change strings) ---
$textValues = @{
    "D2" = "29.388.00"
    "E2" = "  +0.22%  "
    "D3" = "1.883.87"
    "E3" = "  +0.40%  "
    "E4" = "  +0.08%  "
    "E5" = "  +0.22%  "
    "E6" = "  +0.08%  "
    "E7" = "  +0.07%  "
    "E8" = "  +3.83%  "
    "E9" = "  +0.66%  "
    "E10" = "  +0.84%  "
    "E11" = "  -1.19%  "
    "D12" = "1.879.65"
    "E12" = "  -0.62%  "
    "E13" = "  +1.50%  "
    "E14" = "  +0.58%  "
    "E15" = "  +1.41%  "
    "E16" = "  +5.01%  "
    "E17" = "  +2.52%  "
    "D18" = "29.406.46"
    "E18" = "  +0.24%  "
    "B19" = "BitcoinCash"
    "C19" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "E19" = "  -0.47%  "
    "B20" = "WrappedliquidstakedEther2.0"
    "C20" = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
    "D20" = "2.144.55"
    "E20" = "  +0.42%  "
    "E21" = "  +0.16%  "
    "E22" = "  +0.14%  "
    "E23" = "  +0.80%  "
    "E24" = "  +0.05%  "
    "E25" = "  -1.25%  "
    "E26" = "  +0.72%  "
    "E29" = "  -0.47%  "
    "E30" = "  +0.29%  "
    "E31" = "  +0.19%  "
    "E32" = "  -5.21%  "
    "E33" = "  +2.42%  "
    "E34" = "  +1.05%  "
    "E35" = "  +0.57%  "
    "E36" = "  +1.35%  "
    "E37" = "  +0.69%  "
    "E38" = "  +0.95%  "
    "D39" = "1.286.93"
    "E39" = "  +9.70%  "
    "E40" = "  +0.82%  "
    "E41" = "  +3.64%  "
    "E42" = "  +2.02%  "
    "E43" = "  +5.20%  "
    "E44" = "  +0.67%  "
    "E46" = "  +5.97%  "
    "D47" = "2.049.07"
    "E47" = "  +1.01%  "
    "E48" = "  -0.16%  "
    "E49" = "  +0.25%  "
    "E50" = "  +1.04%  "
    "E51" = "  +2.18%  "
}
foreach ($addr in $textValues.Keys) {
    $ws.Range($addr).Value = $textValues[$addr]
}

# --- Cells whose new value looks like a plain number (e.g. "242.48") but must stay as TEXT,
#     matching the original sheet where prices are stored as inline strings, not numbers. ---
$numericLookingValues = @{
    "D5" = "0.7128"
    "D6" = "242.48"
    "D8" = "0.08045"
    "D9" = "0.3130"
    "D10" = "25.26"
    "D11" = "0.08363"
    "D13" = "0.7217"
    "D14" = "5.244"
    "D15" = "92.72"
    "D16" = "6.309"
    "D17" = "0.000008488"
    "D19" = "241.31"
    "D21" = "13.27"
    "D23" = "7.899"
    "D25" = "0.1588"
    "D26" = "163.79"
    "D27" = "9.082"
    "D28" = "18.56"
    "D30" = "4.418"
    "D31" = "4.337"
    "D32" = "1.211"
    "D36" = "0.7496"
    "D37" = "2.703"
    "D38" = "0.01887"
    "D40" = "2.741"
    "D41" = "6.615"
    "D42" = "0.9047"
    "D43" = "111.72"
    "D44" = "73.53"
    "D48" = "1.810"
    "D49" = "0.5217"
    "D50" = "9.499"
    "D51" = "0.4404"
}
foreach ($addr in $numericLookingValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingValues[$addr]
    $cell.Style = "Normal"
}
